# Applies the "Penalty Reward System" (unfinished) edit:
#  - Forecast Comparison sheet: shift Week_Start_Date (col B) one week later
#    and set MyForecast (col D) to 1 for every data row (2-17)
#  - Summary sheet: update several derived metrics in column B

$wb = $excel.ActiveWorkbook

# ---- Sheet: Forecast Comparison ----
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$weekDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

for ($i = 0; $i -lt $weekDates.Length; $i++) {
    $row = $i + 2

    # Write the date as plain text (not an Excel date serial number)
    $cellB = $wsForecast.Range("B$row")
    $cellB.NumberFormat = "@"
    $cellB.Value = $weekDates[$i]
    $cellB.Style = "Normal"

    # MyForecast column becomes 1 for every row
    $wsForecast.Range("D$row").Value = 1
}

# ---- Sheet: Summary ----
$wsSummary = $wb.Worksheets.Item("Summary")

function Set-TextValue($range, $text) {
    # Force the cell to hold a plain text value, matching the source
    # workbook's inline-string cells (avoids Excel auto-converting
    # numeric-looking or date-looking strings into numbers/dates).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $wsSummary.Range("B2")  "2023-02-05 to 2025-01-05"
Set-TextValue $wsSummary.Range("B6")  "8"
Set-TextValue $wsSummary.Range("B9")  "14"
Set-TextValue $wsSummary.Range("B10") "7"
Set-TextValue $wsSummary.Range("B11") "3"
Set-TextValue $wsSummary.Range("B12") "1"
Set-TextValue $wsSummary.Range("B13") "2025-03-23"
Set-TextValue $wsSummary.Range("B14") "1"
Set-TextValue $wsSummary.Range("B15") "2025-04-27"
